# Weekly update: a new "Ajo" price observation (2022-02-25) was recorded at
# the Terminal Hortofrutícola Agro Chillán market. Insert it as the new
# first data row (row 171) for this subset, pushing the existing history
# (previously rows 171-199) down by one row (to 172-200).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 171 - this shifts rows 171:199
# down to 172:200 and keeps the date column's number format (style index 2).
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row with the latest weekly observation.
$ws.Cells.Item(171, 1).Value  = 7
$ws.Cells.Item(171, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(171, 3).Value  = "Ñuble"
$ws.Cells.Item(171, 4).Value  = 44617
$ws.Cells.Item(171, 5).Value  = 16
$ws.Cells.Item(171, 6).Value  = 100112003
$ws.Cells.Item(171, 7).Value  = "Ajo"
$ws.Cells.Item(171, 8).Value  = "Chino"
$ws.Cells.Item(171, 9).Value  = "Primera"
$ws.Cells.Item(171, 10).Value = 60
$ws.Cells.Item(171, 11).Value = 19000
$ws.Cells.Item(171, 12).Value = 20000
$ws.Cells.Item(171, 13).Value = 19500
$ws.Cells.Item(171, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(171, 15).Value = "China"
$ws.Cells.Item(171, 16).Value = 1950
$ws.Cells.Item(171, 17).Value = 10
$ws.Cells.Item(171, 18).Value = "Hortaliza"
